$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Year -> GDP per Capita value (as text, matching the source file's convention
# of storing the Data column values as strings) for 1950-2016.
$data = @{
    1950 = "756"
    1951 = "776"
    1952 = "797"
    1953 = "816"
    1954 = "838"
    1955 = "859"
    1956 = "883"
    1957 = "905"
    1958 = "929"
    1959 = "953"
    1960 = "971"
    1961 = "1004"
    1962 = "1060"
    1963 = "1039"
    1964 = "1055"
    1965 = "1087"
    1966 = "1082"
    1967 = "1165"
    1968 = "1189"
    1969 = "1199"
    1970 = "1188"
    1971 = "1219"
    1972 = "1253"
    1973 = "1172"
    1974 = "1132"
    1975 = "1068"
    1976 = "1039"
    1977 = "1090"
    1978 = "1189"
    1979 = "1173"
    1980 = "1165"
    1981 = "1186"
    1982 = "1183"
    1983 = "1140"
    1984 = "1129"
    1985 = "1243"
    1986 = "1400"
    1987 = "1336"
    1988 = "1403"
    1989 = "1342"
    1990 = "1286"
    1991 = "1331.2303122712"
    1992 = "1264.99198183127"
    1993 = "1240.14216512419"
    1994 = "1189.12132354216"
    1995 = "1191.30736189223"
    1996 = "1255.77492223312"
    1997 = "1268.26320224702"
    1998 = "1293.57340354757"
    1999 = "1318.71976131389"
    2000 = "1275.48190817945"
    2001 = "1283.38218124855"
    2002 = "1255.55709087869"
    2003 = "1277.58086076763"
    2004 = "1268.68662703978"
    2005 = "1309.98364411208"
    2006 = "1322.63719662771"
    2007 = "1308.52031683132"
    2008 = "1315.45371606898"
    2009 = "1287.36782837249"
    2010 = "1326.19751726465"
    2011 = "1344"
    2012 = "1388"
    2013 = "1422"
    2014 = "1437"
    2015 = "1450"
    2016 = "1483"
}

# Country code / name / indicator are constant for every row of this sheet.
$countryCode = 854
$countryName = "Burkina Faso"
$indicator = "GDP per Capita"

# Force the Data column to store these numeric-looking values as text
# (shared strings) rather than numbers, matching the existing rows.
$ws.Range("E2:E68").NumberFormat = "@"

$row = 2
foreach ($year in 1950..2010) {
    $ws.Cells.Item($row, 5).Value = $data[$year]
    $row++
}

# Append the new rows for 2011-2016.
foreach ($year in 2011..2016) {
    $ws.Cells.Item($row, 1).Value = $countryCode
    $ws.Cells.Item($row, 2).Value = $countryName
    $ws.Cells.Item($row, 3).Value = $indicator
    $ws.Cells.Item($row, 4).Value = $year
    $ws.Cells.Item($row, 5).Value = $data[$year]
    $row++
}

# Drop the temporary text format so styling stays as close as possible to
# the original (values keep their text/shared-string type).
$ws.Range("E2:E68").ClearFormats()
